# Add a new "before simulation" header + data row block above the
# existing table (mirrors the header row at row 6, with relax/Hall values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 2) - same labels as the header row at row 6
$ws.Range("B2").Value = "#"
$ws.Range("C2").Value = "gseed"
$ws.Range("D2").Value = "sseed"
$ws.Range("E2").Value = "Nc"
$ws.Range("F2").Value = "Nspcm1"
$ws.Range("G2").Value = "Nspcm2"
$ws.Range("H2").Value = "no sim"
$ws.Range("I2").Value = "act relax"
$ws.Range("J2").Value = "type"
$ws.Range("K2").Value = "kernel"
$ws.Range("L2").Value = "with relax-term-tests"

# Data row (row 3)
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 1001
$ws.Range("D3").Value = 1001
$ws.Range("E3").Value = 75
$ws.Range("F3").Value = 26
$ws.Range("G3").Value = 38
$ws.Range("H3").Value = 50
$ws.Range("I3").Value = "y"
$ws.Range("J3").Value = "AS"
$ws.Range("K3").Value = "Hall"
$ws.Range("L3").Value = "yes"

# Update the active selection to B5 (as in the target file)
$ws.Range("B5").Select()
